$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2023-2024")

# Row 21: new journal entry for the session of 2024-03-04 (TP A2-4-D)
$ws.Range("A21").Value = 45355
$ws.Range("B21").Value = "FSIL"
$ws.Range("C21").Value = "TP"
$ws.Range("E21").Value = "x"

$descr = @"
Suite et fin Elaastic Yasmf 1.6 : All Users (5, 6, 7 + PHPStan et tests PHPUnit sans et avec Coverage lancés).
"@
$ws.Range("G21").Value = $descr

$comment = @"
J'ai insisté sur le fait que l'approche transaction+rollback était indispensable sur les tests d'intégration modififant une BD => à appliquer en SAÉ.
On a pas eu le temps de rebasculer sur question score.
Question : est-ce qu'ils utilisent bien les commits "fix #..." ? J'ai jeté un coup d'oeil au board Git4school et ai vu qu'un seul commit de fix #2...
"@
$ws.Range("I21").Value = $comment

# Update view state to reflect where the author was working
$ws.Range("I22").Select()
